$d = $word.ActiveDocument

# Change 1: insert "production, " after "natural gas " and before "transmission, processing"
# (entirely inside the run that precedes the lastRenderedPageBreak, so it's safe)
$d.Content.Find.Execute(
    "total natural gas transmission,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "total natural gas production, transmission,",
    2)

# Change 2: remove the trailing "are provided by " from the end of that same run
# (still entirely before the page break run)
$d.Content.Find.Execute(
    "North American domain are provided by ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "North American domain ",
    2)

# Change 3: prepend "are provided by " to the start of the run that follows the page
# break (entirely inside that run, so the lastRenderedPageBreak element is untouched)
$d.Content.Find.Execute(
    "the EDGAR v4.3.2 (?) global emission inventory",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "are provided by the EDGAR v4.3.2 (?) global emission inventory",
    2)

# Change 4: merge math runs "= -8." + "18" + " + 0.44" into a single run "= -8.18 + 0.44"
$d.Content.Find.Execute(
    " = -8.1",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    0)
